# CH-154 Custom Index Column
# "Finished my single function"
#
# - Duplicate the "EDA" sheet into a new sheet named "MySingleFunction"
#   (placed after "EDA", i.e. as the last sheet).
# - Add a LET/SCAN based "index" helper column (I17:I27) and a check
#   column comparing it against the existing index (J17:J27) on the
#   new sheet.
# - Update the filter-database defined name for the new sheet.
# - Update sheet selections: EDA is no longer the active tab (selection
#   moves to G36); the new MySingleFunction sheet becomes the active tab
#   (selection N23).

$wb = $excel.ActiveWorkbook

$eda = $wb.Worksheets.Item("EDA")

# Duplicate "EDA" -> new sheet placed right after it (i.e. at the end).
$eda.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$new = $wb.Worksheets.Item($wb.Worksheets.Count)
$new.Name = "MySingleFunction"

# Hidden filter-database range name for the new sheet (mirrors the ones
# that already exist for "EDA" and "Original").
$new.Names.Add("_xlnm._FilterDatabase", "=MySingleFunction!`$E`$2:`$G`$15")

# New helper columns on the copied sheet: a LET/SCAN computed "index"
# (spills I17:I27) and a boolean check against the original index
# column G (spills J17:J27). Entered as (legacy CSE) array formulas so
# they spill exactly like the source workbook's dynamic-array formulas.
$new.Range("I17:I27").FormulaArray = "=LET(d, B3:B13, ds, DROP(VSTACK(0,d),-1), SCAN(0,d-ds,LAMBDA(a,v,IF(v=1,a,a+1))))"
$new.Range("J17:J27").FormulaArray = "=ANCHORARRAY(I17)=G3:G13"

# Selection / active-tab bookkeeping: EDA is no longer selected, the new
# sheet is now the active tab.
$eda.Activate()
$eda.Range("G36").Select()

$new.Activate()
$new.Range("N23").Select()
